# Atualização da tela de matricula e correção de bugs
# Appends attendance records (rows 170-183) to the "Sheet" worksheet,
# mirroring the structure of the existing rows (columns A-F as literal
# text, even when the text looks like a number or a date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($Worksheet, $Address, $Text) {
    # Force the cell to be stored as literal text (matches the source
    # workbook, where every cell - including ids like "44554" and dates
    # like "2024-10-10" - is plain text rather than a number/date).
    $cell = $Worksheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
}

$rows = @(
    @{ Row = 170; A = "44554";  B = "1"; C = "2024-10-10"; D = "P" },
    @{ Row = 171; A = "323443"; B = "3"; C = "2024-10-10"; D = "P" },
    @{ Row = 172; A = "44554";  B = "1"; C = "2024-10-11"; D = "A"; E = "Cirio"; F = "Cirio " },
    @{ Row = 173; A = "323443"; B = "3"; C = "2024-10-11"; D = "A"; E = "Cirio"; F = "Cirio " },
    @{ Row = 174; A = "44554";  B = "1"; C = "2024-10-14"; D = "A"; E = "teste"; F = "teste" },
    @{ Row = 175; A = "323443"; B = "3"; C = "2024-10-14"; D = "A"; E = "teste"; F = "teste" },
    @{ Row = 176; A = "44554";  B = "1"; C = "2024-08-15"; D = "A"; E = "niver nanzao"; F = "niver nanzao" },
    @{ Row = 177; A = "323443"; B = "3"; C = "2024-08-15"; D = "A"; E = "niver nanzao"; F = "niver nanzao" },
    @{ Row = 178; A = "44554";  B = "1"; C = "2024-08-16"; D = "P" },
    @{ Row = 179; A = "323443"; B = "3"; C = "2024-08-16"; D = "P" },
    @{ Row = 180; A = "44554";  B = "1"; C = "2024-08-19"; D = "P" },
    @{ Row = 181; A = "323443"; B = "3"; C = "2024-08-19"; D = "P" },
    @{ Row = 182; A = "44554";  B = "1"; C = "2024-08-20"; D = "A" },
    @{ Row = 183; A = "323443"; B = "3"; C = "2024-08-20"; D = "A" }
)

foreach ($r in $rows) {
    $n = $r.Row
    Set-TextCell $ws "A$n" $r.A
    Set-TextCell $ws "B$n" $r.B
    Set-TextCell $ws "C$n" $r.C
    $ws.Range("D$n").Value = $r.D
    if ($r.ContainsKey("E")) {
        $ws.Range("E$n").Value = $r.E
    }
    if ($r.ContainsKey("F")) {
        $ws.Range("F$n").Value = $r.F
    }
}
